$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 7207.4287
$ws.Range("I2").Value = 7207.4287
$ws.Range("K2").Value = 7207.4287
$ws.Range("M2").Value = -7094.4287
$ws.Range("H17").Value = 1277.5
$ws.Range("I17").Value = 700
$ws.Range("J17").Value = 1999.375
$ws.Range("K17").Value = 2100
$ws.Range("L17").Value = 5998.125
$ws.Range("M17").Value = -1932
$ws.Range("N17").Value = -6334.125
$ws.Range("H33").Value = 101.416664
$ws.Range("I33").Value = 102.09091
$ws.Range("K33").Value = 102.09091
$ws.Range("M33").Value = 126.90909
$ws.Range("H40").Value = 6699.8335
$ws.Range("J40").Value = 8689.200000000001
$ws.Range("L40").Value = 8689.200000000001
$ws.Range("N40").Value = -9039.200000000001
$ws.Range("H51").Value = 10000
$ws.Range("I51").Value = 0
$ws.Range("K51").Value = 0
$ws.Range("M51").ClearContents()
$ws.Range("H137").Value = 2763.65
$ws.Range("I137").Value = 1834.25
$ws.Range("J137").Value = 3383.25
$ws.Range("K137").Value = 5502.75
$ws.Range("L137").Value = 10149.75
$ws.Range("M137").Value = -2952.75
$ws.Range("N137").Value = -15249.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3848323.2
$ws.Range("J32").Value = 20009028
$ws.Range("L32").Value = 20009028
$ws.Range("N32").Value = -20009602
$ws.Range("H61").Value = 2642.3635
$ws.Range("I61").Value = 2717.3
$ws.Range("J61").Value = 1893
$ws.Range("K61").Value = 2717.3
$ws.Range("L61").Value = 1893
$ws.Range("M61").Value = -2505.3
$ws.Range("N61").Value = -2317
$ws.Range("H69").Value = 249999.5
$ws.Range("J69").Value = 249999.5
$ws.Range("L69").Value = 249999.5
$ws.Range("N69").Value = -251497.5
$ws.Range("H72").Value = 249999.5
$ws.Range("J72").Value = 249999.5
$ws.Range("L72").Value = 749998.5
$ws.Range("N72").Value = -757486.5
$ws.Range("H88").Value = 676.875
$ws.Range("I88").Value = 420
$ws.Range("J88").Value = 933.75
$ws.Range("K88").Value = 420
$ws.Range("L88").Value = 933.75
$ws.Range("M88").Value = -14
$ws.Range("N88").Value = -1745.75
$ws.Range("H91").Value = 676.875
$ws.Range("I91").Value = 420
$ws.Range("J91").Value = 933.75
$ws.Range("K91").Value = 420
$ws.Range("L91").Value = 933.75
$ws.Range("M91").Value = 984
$ws.Range("N91").Value = -3741.75
$ws.Range("H136").Value = 2642.3635
$ws.Range("I136").Value = 2717.3
$ws.Range("J136").Value = 1893
$ws.Range("K136").Value = 8151.900000000001
$ws.Range("L136").Value = 5679
$ws.Range("M136").Value = -5601.900000000001
$ws.Range("N136").Value = -10779

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3144.625
$ws.Range("I20").Value = 3426.8572
$ws.Range("K20").Value = 3426.8572
$ws.Range("M20").Value = -3179.8572
$ws.Range("H107").Value = 2274.6667
$ws.Range("I107").Value = 2217.4285
$ws.Range("J107").Value = 2475
$ws.Range("K107").Value = 2217.4285
$ws.Range("L107").Value = 2475
$ws.Range("M107").Value = -297.4285
$ws.Range("N107").Value = -6315

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6588.1514
$ws.Range("I31").Value = 4229.857
$ws.Range("K31").Value = 4229.857
$ws.Range("M31").Value = -3934.857
$ws.Range("H32").Value = 2235.875
$ws.Range("I32").Value = 2235.875
$ws.Range("K32").Value = 2235.875
$ws.Range("M32").Value = -1919.875
$ws.Range("H34").Value = 6588.1514
$ws.Range("I34").Value = 4229.857
$ws.Range("K34").Value = 4229.857
$ws.Range("M34").Value = -4027.857
$ws.Range("H38").Value = 4333
$ws.Range("I38").Value = 3499.5
$ws.Range("J38").Value = 6000
$ws.Range("K38").Value = 3499.5
$ws.Range("L38").Value = 6000
$ws.Range("M38").Value = -3122.5
$ws.Range("N38").Value = -6754
$ws.Range("H44").Value = 18800
$ws.Range("J44").Value = 18800
$ws.Range("L44").Value = 18800
$ws.Range("N44").Value = -19684
$ws.Range("H46").Value = 4333
$ws.Range("I46").Value = 3499.5
$ws.Range("J46").Value = 6000
$ws.Range("K46").Value = 3499.5
$ws.Range("L46").Value = 6000
$ws.Range("M46").Value = -3288.5
$ws.Range("N46").Value = -6422
$ws.Range("H99").Value = 2745.8572
$ws.Range("I99").Value = 1741.4
$ws.Range("K99").Value = 1741.4
$ws.Range("M99").Value = -243.4000000000001
$ws.Range("H106").Value = 72499.5
$ws.Range("J106").Value = 72499.5
$ws.Range("L106").Value = 72499.5
$ws.Range("N106").Value = -75023.5
$ws.Range("H126").Value = 2745.8572
$ws.Range("I126").Value = 1741.4
$ws.Range("K126").Value = 5224.200000000001
$ws.Range("M126").Value = -2754.200000000001
$ws.Range("H134").Value = 2686.25
$ws.Range("I134").Value = 929.5714
$ws.Range("K134").Value = 2788.7142
$ws.Range("M134").Value = -253.7142000000003

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 66603.625
$ws.Range("I132").Value = 80512.16
$ws.Range("J132").Value = 6333.3335
$ws.Range("K132").Value = 241536.48
$ws.Range("L132").Value = 19000.0005
$ws.Range("M132").Value = -239006.48
$ws.Range("N132").Value = -24060.0005

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1590
$ws.Range("I16").Value = 1590
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 1590
$ws.Range("L16").Value = 0
$ws.Range("M16").ClearContents()
$ws.Range("N16").Value = -1420
$ws.Range("H22").Value = 1465.3529
$ws.Range("I22").Value = 1114
$ws.Range("K22").Value = 1114
$ws.Range("M22").Value = -819
$ws.Range("H27").Value = 1465.3529
$ws.Range("I27").Value = 1114
$ws.Range("K27").Value = 1114
$ws.Range("M27").Value = -1007
$ws.Range("H40").Value = 10004
$ws.Range("I40").Value = 10004
$ws.Range("K40").Value = 10004
$ws.Range("M40").Value = -9868
$ws.Range("H46").Value = 767.3333
$ws.Range("I46").Value = 550.3333
$ws.Range("J46").Value = 875.8333
$ws.Range("K46").Value = 550.3333
$ws.Range("L46").Value = 875.8333
$ws.Range("M46").Value = -362.3333
$ws.Range("N46").Value = -1251.8333
$ws.Range("H101").Value = 14405.2
$ws.Range("J101").Value = 14858.111
$ws.Range("L101").Value = 14858.111
$ws.Range("N101").Value = -21348.111
$ws.Range("H122").Value = 10004
$ws.Range("I122").Value = 10004
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 30012
$ws.Range("L122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("N122").Value = -27562
$ws.Range("H136").Value = 0
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("L136").ClearContents()
$ws.Range("M136").ClearContents()
$ws.Range("N136").Value = 0

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 7704.727
$ws.Range("I62").Value = 3150.6
$ws.Range("K62").Value = 3150.6
$ws.Range("M62").Value = -2526.6
$ws.Range("H65").Value = 7704.727
$ws.Range("I65").Value = 3150.6
$ws.Range("K65").Value = 15753
$ws.Range("M65").Value = -12633
$ws.Range("H81").Value = 749
$ws.Range("I81").Value = 749
$ws.Range("K81").Value = 1498
$ws.Range("M81").Value = -437
$ws.Range("H84").Value = 749
$ws.Range("I84").Value = 749
$ws.Range("K84").Value = 7490
$ws.Range("M84").Value = -2186
$ws.Range("H122").Value = 3979
$ws.Range("I122").Value = 2583.3333
$ws.Range("K122").Value = 7749.999899999999
$ws.Range("M122").Value = -5299.999899999999
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").ClearContents()
$ws.Range("N123").Value = 0
$ws.Range("H132").Value = 5002
$ws.Range("I132").Value = 5002
$ws.Range("K132").Value = 15006
$ws.Range("M132").Value = -12476
